$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 2466
$ws.Cells.Item(5, 6).Value = 1752
$ws.Cells.Item(6, 6).Value = 108
$ws.Cells.Item(7, 6).Value = 323
$ws.Cells.Item(9, 6).Value = 3568
$ws.Cells.Item(10, 6).Value = 1190
$ws.Cells.Item(11, 6).Value = 1584
$ws.Cells.Item(15, 6).Value = 1346
$ws.Cells.Item(16, 6).Value = 1796
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(19, 6).Value = 468
$ws.Cells.Item(21, 6).Value = 21
$ws.Cells.Item(22, 6).Value = 2
$ws.Cells.Item(24, 6).Value = 2320
$ws.Cells.Item(25, 6).Value = 151
$ws.Cells.Item(27, 6).Value = 4336
$ws.Cells.Item(28, 6).Value = 55
$ws.Cells.Item(30, 6).Value = 2693
$ws.Cells.Item(33, 6).Value = 1222
$ws.Cells.Item(34, 6).Value = 935
$ws.Cells.Item(35, 6).Value = 6

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(14, 6).Value = 39
$ws.Cells.Item(22, 6).Value = 135
$ws.Cells.Item(35, 6).Value = 442
$ws.Cells.Item(39, 6).Value = 29
$ws.Cells.Item(42, 6).Value = 86

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 2543
$ws.Cells.Item(5, 6).Value = 9580
$ws.Cells.Item(9, 6).Value = 393
$ws.Cells.Item(10, 6).Value = 2997
$ws.Cells.Item(11, 6).Value = 511
$ws.Cells.Item(12, 6).Value = 819
$ws.Cells.Item(13, 6).Value = 230
$ws.Cells.Item(14, 6).Value = 260

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 2543
$ws.Cells.Item(5, 6).Value = 2466
$ws.Cells.Item(8, 6).Value = 393
$ws.Cells.Item(9, 6).Value = 2997
$ws.Cells.Item(10, 6).Value = 511
$ws.Cells.Item(11, 6).Value = 819
$ws.Cells.Item(12, 6).Value = 230
$ws.Cells.Item(15, 6).Value = 1752
$ws.Cells.Item(16, 6).Value = 323
$ws.Cells.Item(23, 6).Value = 39
$ws.Cells.Item(26, 6).Value = 1796
$ws.Cells.Item(29, 6).Value = 135
$ws.Cells.Item(30, 6).Value = 135
$ws.Cells.Item(31, 6).Value = 21
$ws.Cells.Item(35, 6).Value = 2320
$ws.Cells.Item(39, 6).Value = 260
$ws.Cells.Item(43, 6).Value = 29
$ws.Cells.Item(46, 6).Value = 86
$ws.Cells.Item(49, 6).Value = 1222
$ws.Cells.Item(50, 6).Value = 935
$ws.Cells.Item(51, 6).Value = 6
